$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

$ws.Range("D2").Value = "/Users/runner/runners/2.263.0/work/_temp/Library/AlpsNMR/dataset-demo/10.zip"
$ws.Range("D3").Value = "/Users/runner/runners/2.263.0/work/_temp/Library/AlpsNMR/dataset-demo/20.zip"
$ws.Range("D4").Value = "/Users/runner/runners/2.263.0/work/_temp/Library/AlpsNMR/dataset-demo/30.zip"
$ws.Range("D5").Value = "/Users/runner/runners/2.263.0/work/_temp/Library/AlpsNMR/dataset-demo/40.zip"
